$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.906.72"
$ws.Range("E2").Value = "  -3.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.49"
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.93"
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4368"
$ws.Range("E7").Value = "  -4.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3704"
$ws.Range("E8").Value = "  -3.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07500"
$ws.Range("E9").Value = "  -2.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9382"
$ws.Range("E10").Value = "  -4.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.36"
$ws.Range("E11").Value = "  -3.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.882.46"
$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.752"
$ws.Range("E13").Value = "  -2.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.457"
$ws.Range("E14").Value = "  -3.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06833"
$ws.Range("E15").Value = "  -2.71%  "

$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.64"
$ws.Range("E17").Value = "  -2.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009095"
$ws.Range("E18").Value = "  -3.77%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.93"
$ws.Range("E20").Value = "  -4.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.899.20"
$ws.Range("E21").Value = "  -3.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.116"
$ws.Range("E22").Value = "  -3.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.10"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.111.34"
$ws.Range("E24").Value = "  -1.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.003"
$ws.Range("E25").Value = "  -4.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.99"
$ws.Range("E26").Value = "  -2.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.494"
$ws.Range("E28").Value = "  -3.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.20"
$ws.Range("E29").Value = "  -3.57%  "

$ws.Range("E30").Value = "  -8.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09035"
$ws.Range("E31").Value = "  -2.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8156"
$ws.Range("E32").Value = "  -5.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.821"
$ws.Range("E33").Value = "  -5.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.175"
$ws.Range("E34").Value = "  -5.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.953"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("E36").Value = "  +0.33%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.120"
$ws.Range("E37").Value = "  -3.17%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05482"
$ws.Range("E38").Value = "  -3.86%  "

$ws.Range("E39").Value = "  -2.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.938"
$ws.Range("E40").Value = "  -0.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5267"
$ws.Range("E41").Value = "  -4.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.035"
$ws.Range("E42").Value = "  -6.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1706"
$ws.Range("E43").Value = "  -2.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.798"
$ws.Range("E44").Value = "  -6.12%  "

$ws.Range("E45").Value = "  -1.70%  "

$ws.Range("E46").Value = "  -4.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.67"
$ws.Range("E47").Value = "  -4.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.07"
$ws.Range("E48").Value = "  -2.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.682"
$ws.Range("E49").Value = "  -5.48%  "

$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.886"
$ws.Range("E51").Value = "  -12.97%  "
